$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells for the "highest_pl_rank_score" / "highest_pl_rank" columns
$ws.Range("E3").Value = "pl_score"
$ws.Range("F3").Value = "pl_rank"

# Update per-player current PL rank score/rank (was highest_pl_rank*, now current pl_*)
# Row 5
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = "Mythic I"

# Row 6
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = ""

# Row 7
$ws.Range("E7").Value = 14
$ws.Range("F7").Value = "Mythic II"

# Row 8
$ws.Range("D8").Value = 46017
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = "Legend I"

# Row 9
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "Mythic III"

# Row 10
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = 42

# Row 11
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = "Gold III"
$ws.Range("G11").Value = 12

# Row 12
$ws.Range("E12").Value = 12
$ws.Range("F12").Value = "Diamond III"

# Row 13
$ws.Range("E13").Value = 11
$ws.Range("F13").Value = "Diamond II"
$ws.Range("I13").Value = 61

# Row 14
$ws.Range("E14").Value = 14
$ws.Range("F14").Value = "Mythic II"

# Row 15
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = "Mythic I"

# Row 16
$ws.Range("E16").Value = 12
$ws.Range("F16").Value = "Diamond III"

# Row 17
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = "Mythic I"

# Row 18
$ws.Range("D18").Value = 42700
$ws.Range("E18").Value = 14
$ws.Range("F18").Value = "Mythic II"

# Row 19
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "Mythic III"

# Row 20
$ws.Range("E20").Value = 14
$ws.Range("F20").Value = "Mythic II"

# Row 21
$ws.Range("D21").Value = 21837
$ws.Range("F21").Value = "Mythic III"

# Row 22
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = "Mythic I"

# Row 23
$ws.Range("E23").Value = 11
$ws.Range("F23").Value = "Diamond II"
$ws.Range("I23").Value = 32

# Row 24
$ws.Range("D24").Value = 31533
$ws.Range("E24").Value = 12
$ws.Range("F24").Value = "Diamond III"

# Row 25
$ws.Range("E25").Value = 12
$ws.Range("F25").Value = "Diamond III"

# Row 26
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = ""

# Row 27
$ws.Range("E27").Value = 11
$ws.Range("F27").Value = "Diamond II"
$ws.Range("I27").Value = 48

# Row 28
$ws.Range("E28").Value = 11
$ws.Range("F28").Value = "Diamond II"
$ws.Range("I28").Value = 40

# Row 30
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = "Diamond I"
$ws.Range("H30").Value = 14

# Row 32
$ws.Range("E32").Value = 12
$ws.Range("F32").Value = "Diamond III"

# Row 33
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = ""
